# Update the three message templates on the "Messages" sheet with
# the new [vb]-delimited formatting, matching the author's edit.

$customerMsg = @"
Dear [Customer's Name],
[vb]
Thank you for your order! We’re excited to let you know that your order is currently being prepared. Delivery typically takes between 1 to 3 days.
[vb]
Order Details >>
[vb]
Order Number: [Order Number] ,
[vb]
[Order items]
[vb]
Delivery Fee: [Delivery Fee]
[vb]
Total Price: [Total Price for All Items]
[vb]
[vb]
If you have any questions or need further assistance, feel free to contact us at [Website Contact Info].
[vb]
Thank you for choosing Bookle Store!
[vb]
Best regards,
[vb]
Customer Care Team
[vb]
Bookle Store
"@

$driverMsg = @"
Dear [Driver's Name],
[vb]
You have a new delivery assignment. Below are the details:
[vb]
Order Details >>
[vb]
Order Number: [Order Number] ,
[vb]
Customer Name: [Customer's Name] ,
[vb]
Customer Address: [Customer's Address] ,
[vb]
Customer Phone Number: [Customer's Phone Number] ,
[vb]
Delivery Timeframe: 1-3 days
Items to Deliver:
[vb]
[Order items]
[vb]
[vb]
Please ensure all items are picked up from the respective publishers and delivered to the customer on time. If you encounter any issues, don’t hesitate to contact us at [Website Contact Info].
[vb]
Thank you for your efforts!
[vb]
Best regards,
[vb]
Customer Support Team
[vb]
Bookle Store
"@

$publisherMsg = @"
Dear [Publisher's Name],
[vb]
We hope this message finds you well. We have a new order that includes items from your inventory. Below are the details:
[vb]
Order Details >>
[vb]
Order Number: [Order Number] ,
[vb]
Customer Name: [Customer's Name] ,
[vb]
Customer Phone Number: [Customer's Phone Number] ,
[vb]
Delivery Timeframe: 1-3 days
Books to Prepare:
[vb]
[Order items]
[vb]
Please ensure the items are ready for pickup by our delivery driver. If you have any questions or need clarification, feel free to reach out to us at [Website Contact Info].
[vb]
[vb]
Thank you for your cooperation!
[vb]
Best regards,
[vb]
Customer Support Team
[vb]
Bookle Store
"@

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Messages")
$ws.Activate()

$ws.Range("A2").Value = $customerMsg
$ws.Range("B2").Value = $driverMsg
$ws.Range("C2").Value = $publisherMsg

# The longer templates push row 2 to Excel's maximum row height.
$ws.Rows.Item(2).RowHeight = 409.6

# Matches the saved selection/scroll state in the authored workbook
# (scrolled down so row 2 is the topmost visible row).
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
